$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each entry: row, column, new text. Addressed by fixed table position so
# that replacement text which happens to equal another cell's old text
# (e.g. "59÷5=11, 4") cannot be re-matched by a later step.
$updates = @(
    @(1, 1, "59÷5=11, 4"),
    @(1, 2, "47÷4=11, 3"),
    @(1, 3, "32÷8=4, 0"),
    @(1, 4, "97÷8=12, 1"),
    @(1, 5, "20÷7=2, 6"),

    @(5, 1, "19÷6=3, 1"),
    @(5, 2, "89÷4=22, 1"),
    @(5, 3, "50÷5=10, 0"),
    @(5, 4, "57÷9=6, 3"),
    @(5, 5, "51÷7=7, 2"),

    @(9, 1, "88÷8=11, 0"),
    @(9, 2, "83÷9=9, 2"),
    @(9, 3, "12÷8=1, 4"),
    @(9, 4, "76÷3=25, 1"),
    @(9, 5, "55÷8=6, 7"),

    @(13, 1, "71÷8=8, 7"),
    @(13, 2, "57÷6=9, 3"),
    @(13, 3, "73÷4=18, 1"),
    @(13, 4, "71÷5=14, 1"),
    @(13, 5, "84÷4=21, 0"),

    @(17, 1, "23÷9=2, 5"),
    @(17, 2, "45÷5=9, 0"),
    @(17, 3, "66÷5=13, 1"),
    @(17, 4, "76÷3=25, 1"),
    @(17, 5, "15÷9=1, 6")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $newText = $u[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}
